$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing row 236 values (open/high/low/close all become 177772000000000)
$ws.Range("C236").Value = 177772000000000
$ws.Range("D236").Value = 177772000000000
$ws.Range("E236").Value = 177772000000000
$ws.Range("F236").Value = 177772000000000

# New rows of data appended to the table
$newRows = @(
    @{ Row = 237; Date = 45108.41666666666; Val = 176788000000000 },
    @{ Row = 238; Date = 45139.41666666666; Val = 174322000000000 },
    @{ Row = 239; Date = 45170.41666666666; Val = 173950000000000 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    # Copy formatting from the row above so the new row matches existing style
    $ws.Range("A$($r - 1):G$($r - 1)").Copy() | Out-Null
    $ws.Range("A$($r):G$($r)").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Range("A$r").Value = $item.Date
    $ws.Range("B$r").Value = "ECONOMICS:IQM2"
    $ws.Range("C$r").Value = $item.Val
    $ws.Range("D$r").Value = $item.Val
    $ws.Range("E$r").Value = $item.Val
    $ws.Range("F$r").Value = $item.Val
    $ws.Range("G$r").Value = 0
}

$excel.CutCopyMode = 0
